# flavi-ncbi-refseqs-side-data.xlsx -- "Refactor, build AA trees and update website"
#
# Updates the vector "clade" (column H) for a block of Mosq1 (G="Mosq1") rows
# on the flavi.txt sheet from the placeholder NULL down to the correct
# mosquito genus (Culex / Aedes), normalizes the formatting on a couple of
# cells whose styling had drifted from the rest of their row-band, and moves
# the frozen-pane scroll position / active selection to where the editor was
# last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("flavi.txt")

# --- Column H ("clade") corrections, rows 29-63 -----------------------------
# row -> new clade value
$cladeUpdates = @{
    29 = "Culex"
    30 = "Culex"
    31 = "Aedes"
    32 = "Culex"
    33 = "Culex"
    34 = "Aedes"
    35 = "Aedes"
    36 = "Aedes"
    37 = "Aedes"
    38 = "Culex"
    39 = "Culex"
    40 = "Culex"
    41 = "Aedes"
    42 = "Culex"
    43 = "Culex"
    44 = "Culex"
    45 = "Aedes"
    46 = "Culex"
    47 = "Culex"
    48 = "Aedes"
    49 = "Culex"
    50 = "Culex"
    51 = "Culex"
    52 = "Culex"
    53 = "Culex"
    54 = "Culex"
    55 = "Culex"
    56 = "Aedes"
    57 = "Culex"
    58 = "Culex"
    59 = "Culex"
    60 = "Culex"
    61 = "Culex"
    62 = "Culex"
    63 = "Culex"
}

foreach ($row in $cladeUpdates.Keys) {
    $ws.Range("H$row").Value = $cladeUpdates[$row]
}

# --- Formatting cleanup -------------------------------------------------
# Row 45's H cell had picked up the neighbouring "left aligned" style (like
# G45) instead of the plain style shared by the rest of the row-band; pull
# the formatting back in line with its siblings (e.g. H44).
$ws.Range("H44").Copy()
$ws.Range("H45").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H45").Value = $cladeUpdates[45]

# Rows 64 and 75 have an I cell whose fill didn't match the rest of the
# row-band (H64/H75) -- bring it into line.
$ws.Range("H64").Copy()
$ws.Range("I64").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H75").Copy()
$ws.Range("I75").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- View state: scroll the frozen pane back to the top and move the ------
# active selection to where editing continued.
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$win.Zoom = 95
[void]$ws.Range("D31").Select()
